$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.981.25"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.83%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.878.63"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.96%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9981"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.09"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -4.65%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9977"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.19%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4993"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -2.62%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2923"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.75%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06619"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.99%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.879.01"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.94%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "16.73"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -3.86%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07246"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.48%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6673"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -3.74%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "86.22"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.88%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.868"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.04%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "29.946.86"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000007908"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.29%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9975"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.29%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.76"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.08%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.120.80"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.89%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9969"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.753"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.643"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.62%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.048"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.49%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "149.20"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.87%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "140.71"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.15%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.09"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.09%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.909"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -5.38%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.390"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.43%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.171"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.71%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08782"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.60%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.944"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.10%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05063"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.47%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7104"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.03%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -4.78%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.85%  "
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "MXToken"
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.692"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -5.05%  "
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01749"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.62%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.180"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -5.71%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9294"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -4.91%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.782"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -5.36%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4250"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.57%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.96"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -4.06%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.453"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -3.44%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1256"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05653"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.54%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "32.40"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -3.23%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.3751"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.19%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.174"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -4.60%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "55.77"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.98%  "
